$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the next empty row after the existing data (row 39 -> 40)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Write the timestamp as plain text (not a date) into column A
$ws.Cells.Item($newRow, 1).Value = "2025-04-29 03:45:14"

# Write the metric value as a number into column B
$ws.Cells.Item($newRow, 2).Value = 106
